$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "1.00", "553.90") are preserved exactly as typed, then restore
# the default "Normal" style so no stray s="n" attribute is left on cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.820.55"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "3.470.66"
$ws.Range("E3").Value = "  +5.54%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "553.90"
$ws.Range("E5").Value = "  +5.74%  "
$ws.Range("D6").Value = "181.25"
$ws.Range("E6").Value = "  +5.69%  "
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +10.43%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("E10").Value = "  +14.55%  "
$ws.Range("D11").Value = "53.99"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").Value = "0.0000273"
$ws.Range("E12").Value = "  +6.85%  "
$ws.Range("D13").Value = "9.28"
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").Value = "4.032.08"
$ws.Range("E14").Value = "  +5.71%  "
$ws.Range("D15").Value = "3.472.31"
$ws.Range("E15").Value = "  +5.71%  "
$ws.Range("D16").Value = "18.53"
$ws.Range("E16").Value = "  +7.30%  "
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "65.882.95"
$ws.Range("E18").Value = "  +3.18%  "
$ws.Range("E19").Value = "  +7.99%  "
$ws.Range("D20").Value = "0.992"
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("D21").Value = "420.40"
$ws.Range("E21").Value = "  +11.37%  "
$ws.Range("E22").Value = "  +10.48%  "
$ws.Range("D23").Value = "85.71"
$ws.Range("E23").Value = "  +5.85%  "
$ws.Range("D24").Value = "4.12"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("D26").Value = "10.82"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "12.28"
$ws.Range("E27").Value = "  +10.24%  "
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "9.03"
$ws.Range("E29").Value = "  +11.90%  "
$ws.Range("D30").Value = "30.10"
$ws.Range("E30").Value = "  +5.07%  "
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").Value = "618.01"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").Value = "11.76"
$ws.Range("E33").Value = "  +5.91%  "
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("D35").Value = "59.97"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("D36").Value = "0.149"
$ws.Range("E36").Value = "  +19.45%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "37.62"
$ws.Range("E38").Value = "  +5.42%  "
$ws.Range("E39").Value = "  +6.13%  "
$ws.Range("D40").Value = "0.382"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("D42").Value = "3.130.51"
$ws.Range("E42").Value = "  +8.10%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +9.47%  "
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "3.29"
$ws.Range("E46").Value = "  +4.35%  "
$ws.Range("E47").Value = "  +4.71%  "
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  +7.00%  "
$ws.Range("D50").Value = "138.62"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "8.40"
$ws.Range("E51").Value = "  +5.74%  "

$ws.Range("D2:D51").Style = "Normal"
